$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1) - row 2 updates
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 973.9537847600009
$ws.Range("E2").Value = 28982.37596598056
$ws.Range("I2").Value = 16175.28135478
$ws.Range("L2").Value = 48524.529503538
$ws.Range("M2").Value = 10590.587968015
$ws.Range("N2").Value = 7155.07579047334
$ws.Range("O2").Value = 6980.325566461758

# Sheet "2030" (sheet2) - row 2 updates
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5712.560177842886
$ws.Range("E2").Value = 56106.05588781912
$ws.Range("I2").Value = 44217.8984721661
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 21984.28023276101
$ws.Range("N2").Value = 10592.74688452318
$ws.Range("O2").Value = 12062.16651258332

# Sheet "2035" (sheet3) - row 2 updates
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15130.68721665935
$ws.Range("O2").Value = 14760.15862166215

# Sheet "2040" (sheet4) - row 2 updates
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15232.90602690694
$ws.Range("O2").Value = 14760.15862166215

# Sheet "2045" (sheet5) - row 2 updates
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15768.41721148085
$ws.Range("O2").Value = 17095.62460801082

# Sheet "2050" (sheet6) - row 2 updates
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15768.41721148085
$ws.Range("O2").Value = 17095.62460801082
